# "finestra incidenza 7gg centrata su ultimo g"
#
# Column B holds daily new-positive counts (rows 2..184, one row per day).
# Columns C/D hold a rolling 7-day-window sum (C) and the same sum scaled
# to cases-per-100000-inhabitants (D).
#
# Previously the 7-day window was CENTERED on the current day (3 days
# before .. 3 days after). This change re-centers the window so it ends ON
# the current/last day of the window (6 days before .. current day), i.e.
# a trailing window: C[r] = SUM(B[r-6 .. r]).
#
# Rows that don't have 6 prior days of history (the first 6 data rows) no
# longer have a defined value and are cleared; rows that now DO have 6
# prior days of history (the last 3 data rows, which previously lacked
# 3 days of "future" data for the centered window) get real values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 184
$windowDays = 7
$population = 686

# Read column B (new positives) into a lookup table keyed by row number.
$newPos = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $newPos[$r] = $ws.Cells.Item($r, 2).Value2
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $windowStart = $r - ($windowDays - 1)

    if ($windowStart -lt $firstRow) {
        # Not enough trailing history yet for a full 7-day window.
        $ws.Cells.Item($r, 3).ClearContents()
        $ws.Cells.Item($r, 4).ClearContents()
    }
    else {
        $sum = 0
        for ($i = $windowStart; $i -le $r; $i++) {
            $sum = $sum + $newPos[$i]
        }

        $perHundredK = $sum * 100000 / $population

        $ws.Cells.Item($r, 3).Value = $sum
        $ws.Cells.Item($r, 4).Value = $perHundredK
    }
}
